{"js": "// Lattice-multiplication worksheet: replace the multiplication problem in\n// each table cell with a new problem, preserving the existing \"AA x BB /\n// digits / ---- / digit| | / digit| |\" layout (5 lines separated by manual\n// line breaks) and the 16pt (sz=32 half-points) run formatting.\n//\n// New content, in row-major order (5 rows x 3 columns):\nconst cellText = [\n  [\"44 x 46\\u000b  4    6\\u000b  ----\\u000b4|    |\\u000b4|    |\", \"14 x 13\\u000b  1    3\\u000b  ----\\u000b1|    |\\u000b4|    |\", \"63 x 84\\u000b  8    4\\u000b  ----\\u000b6|    |\\u000b3|    |\"],\n  [\"13 x 76\\u000b  7    6\\u000b  ----\\u000b1|    |\\u000b3|    |\", \"28 x 90\\u000b  9    0\\u000b  ----\\u000b2|    |\\u000b8|    |\", \"30 x 42\\u000b  4    2\\u000b  ----\\u000b3|    |\\u000b0|    |\"],\n  [\"70 x 69\\u000b  6    9\\u000b  ----\\u000b7|    |\\u000b0|    |\", \"92 x 72\\u000b  7    2\\u000b  ----\\u000b9|    |\\u000b2|    |\", \"65 x 27\\u000b  2    7\\u000b  ----\\u000b6|    |\\u000b5|    |\"],\n  [\"18 x 59\\u000b  5    9\\u000b  ----\\u000b1|    |\\u000b8|    |\", \"56 x 22\\u000b  2    2\\u000b  ----\\u000b5|    |\\u000b6|    |\", \"95 x 88\\u000b  8    8\\u000b  ----\\u000b9|    |\\u000b5|    |\"],\n  [\"16 x 55\\u000b  5    5\\u000b  ----\\u000b1|    |\\u000b6|    |\", \"10 x 32\\u000b  3    2\\u000b  ----\\u000b1|    |\\u000b0|    |\", \"96 x 39\\u000b  3    9\\u000b  ----\\u000b9|    |\\u000b6|    |\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.rows.load(\"items\");\nawait context.sync();\n\nconst rowCount = table.rows.items.length;\n\nfor (let r = 0; r < rowCount; r++) {\n  const row = table.rows.items[r];\n  row.cells.load(\"items\");\n  await context.sync();\n\n  const colCount = row.cells.items.length;\n  for (let c = 0; c < colCount; c++) {\n    if (r >= cellText.length || c >= cellText[r].length) continue;\n\n    const cell = row.cells.items[c];\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n\n    // Replace the paragraph's text while keeping the cell's paragraph/table\n    // structure intact.\n    range.insertText(cellText[r][c], \"Replace\");\n    await context.sync();\n\n    // Re-apply the 16pt run formatting (\"sz\" = 32 half-points) that the\n    // original problem text used, since Replace can drop run properties.\n    range.font.size = 16;\n    await context.sync();\n  }\n}\n", "ps1": "# Lattice-multiplication worksheet: replace the multiplication problem in\n# each table cell with a new problem, preserving the existing \"AA x BB /\n# digits / ---- / digit| | / digit| |\" layout (5 lines separated by manual\n# line breaks, char 11 / `v) and the 16pt run formatting already applied to\n# the cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# New content, in row-major order (5 rows x 3 columns). `v is the\n# PowerShell escape for a vertical-tab (chr 11), which Word renders as a\n# manual line break (<w:br/>) inside a table-cell paragraph.\n$cellText = @(\n  @(\"44 x 46`v  4    6`v  ----`v4|    |`v4|    |\", \"14 x 13`v  1    3`v  ----`v1|    |`v4|    |\", \"63 x 84`v  8    4`v  ----`v6|    |`v3|    |\"),\n  @(\"13 x 76`v  7    6`v  ----`v1|    |`v3|    |\", \"28 x 90`v  9    0`v  ----`v2|    |`v8|    |\", \"30 x 42`v  4    2`v  ----`v3|    |`v0|    |\"),\n  @(\"70 x 69`v  6    9`v  ----`v7|    |`v0|    |\", \"92 x 72`v  7    2`v  ----`v9|    |`v2|    |\", \"65 x 27`v  2    7`v  ----`v6|    |`v5|    |\"),\n  @(\"18 x 59`v  5    9`v  ----`v1|    |`v8|    |\", \"56 x 22`v  2    2`v  ----`v5|    |`v6|    |\", \"95 x 88`v  8    8`v  ----`v9|    |`v5|    |\"),\n  @(\"16 x 55`v  5    5`v  ----`v1|    |`v6|    |\", \"10 x 32`v  3    2`v  ----`v1|    |`v0|    |\", \"96 x 39`v  3    9`v  ----`v9|    |`v6|    |\")\n)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        # Assigning to Range.Text replaces the cell's text in place while\n        # keeping the run's existing character formatting (sz=32).\n        $cell.Range.Text = $cellText[$r - 1][$c - 1]\n    }\n}\n"}
